$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename two field labels in column A (column B "Type" values are unchanged):
#   row 4:  family_id -> referral_id
#   row 17: assembly  -> genome_build
$ws.Range("A4").Value = "referral_id"
$ws.Range("A17").Value = "genome_build"

# Move the active selection to B17, matching the saved workbook state.
$ws.Range("B17").Select() | Out-Null
